$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"22.95531766666667"
$ws.Range("H2").Value = [double]"68.865953"
$ws.Range("I2").Value = [double]"0.1720020945576478"
$ws.Range("J2").Value = [double]"0.1720020945576478"
$ws.Range("M2").Value = [double]"0.1825283333333333"
$ws.Range("N2").Value = [double]"0.547585"
$ws.Range("O2").Value = [double]"0.001028331058213739"
$ws.Range("P2").Value = [double]"0.001028331058213739"
$ws.Range("Q2").Value = [double]"4.189995874833889"
$ws.Range("R2").Value = [double]"37.709962873505"
$ws.Range("S2").Value = [double]"0.0001768750959114455"
$ws.Range("T2").Value = [double]"0.0001768750959114456"

$ws.Range("G3").Value = [double]"22.95531766666667"
$ws.Range("H3").Value = [double]"68.865953"
$ws.Range("I3").Value = [double]"0.1720020945576478"
$ws.Range("J3").Value = [double]"0.1720020945576478"
$ws.Range("O3").Value = [double]"0.0001759459539160193"
$ws.Range("P3").Value = [double]"0.0001759459539160193"
$ws.Range("Q3").Value = [double]"0.7169022225025556"
$ws.Range("R3").Value = [double]"6.452120002523"
$ws.Range("S3").Value = [double]"3.026307260249869E-05"
$ws.Range("T3").Value = [double]"3.026307260249869E-05"

$ws.Range("G4").Value = [double]"22.95531766666667"
$ws.Range("H4").Value = [double]"68.865953"
$ws.Range("I4").Value = [double]"0.1720020945576478"
$ws.Range("J4").Value = [double]"0.1720020945576478"
$ws.Range("M4").Value = [double]"103.239782"
$ws.Range("N4").Value = [double]"309.719346"
$ws.Range("O4").Value = [double]"0.5816339432625932"
$ws.Range("P4").Value = [double]"0.5816339432625932"
$ws.Range("Q4").Value = [double]"2369.901991647415"
$ws.Range("R4").Value = [double]"21329.11792482674"
$ws.Range("S4").Value = [double]"0.1000422565069901"
$ws.Range("T4").Value = [double]"0.1000422565069901"

$ws.Range("G5").Value = [double]"22.95531766666667"
$ws.Range("H5").Value = [double]"68.865953"
$ws.Range("I5").Value = [double]"0.1720020945576478"
$ws.Range("J5").Value = [double]"0.1720020945576478"
$ws.Range("M5").Value = [double]"0.04852733333333333"
$ws.Range("N5").Value = [double]"0.145582"
$ws.Range("O5").Value = [double]"0.0002733940705404138"
$ws.Range("P5").Value = [double]"0.0002733940705404139"
$ws.Range("Q5").Value = [double]"1.113960352182889"
$ws.Range("R5").Value = [double]"10.025643169646"
$ws.Range("S5").Value = [double]"4.702435277259251E-05"
$ws.Range("T5").Value = [double]"4.702435277259251E-05"

$ws.Range("G6").Value = [double]"22.95531766666667"
$ws.Range("H6").Value = [double]"68.865953"
$ws.Range("I6").Value = [double]"0.1720020945576478"
$ws.Range("J6").Value = [double]"0.1720020945576478"
$ws.Range("M6").Value = [double]"73.99751433333334"
$ws.Range("N6").Value = [double]"221.992543"
$ws.Range("O6").Value = [double]"0.4168883856547366"
$ws.Range("P6").Value = [double]"0.4168883856547366"
$ws.Range("Q6").Value = [double]"1698.636448065387"
$ws.Range("R6").Value = [double]"15287.72803258848"
$ws.Range("S6").Value = [double]"0.07170567552937117"
$ws.Range("T6").Value = [double]"0.07170567552937117"

$ws.Range("I7").Value = [double]"0.4661646602805707"
$ws.Range("J7").Value = [double]"0.4661646602805707"
$ws.Range("M7").Value = [double]"0.1825283333333333"
$ws.Range("N7").Value = [double]"0.547585"
$ws.Range("O7").Value = [double]"0.001028331058213739"
$ws.Range("P7").Value = [double]"0.001028331058213739"
$ws.Range("Q7").Value = [double]"11.35583847738722"
$ws.Range("R7").Value = [double]"102.202546296485"
$ws.Range("S7").Value = [double]"0.0004793715984081672"
$ws.Range("T7").Value = [double]"0.0004793715984081673"

$ws.Range("I8").Value = [double]"0.4661646602805707"
$ws.Range("J8").Value = [double]"0.4661646602805707"
$ws.Range("O8").Value = [double]"0.0001759459539160193"
$ws.Range("P8").Value = [double]"0.0001759459539160193"
$ws.Range("S8").Value = [double]"8.201978583500207E-05"
$ws.Range("T8").Value = [double]"8.201978583500207E-05"

$ws.Range("I9").Value = [double]"0.4661646602805707"
$ws.Range("J9").Value = [double]"0.4661646602805707"
$ws.Range("M9").Value = [double]"103.239782"
$ws.Range("N9").Value = [double]"309.719346"
$ws.Range("O9").Value = [double]"0.5816339432625932"
$ws.Range("P9").Value = [double]"0.5816339432625932"
$ws.Range("Q9").Value = [double]"6422.971532269887"
$ws.Range("R9").Value = [double]"57806.74379042898"
$ws.Range("S9").Value = [double]"0.2711371895686555"
$ws.Range("T9").Value = [double]"0.2711371895686555"

$ws.Range("I10").Value = [double]"0.4661646602805707"
$ws.Range("J10").Value = [double]"0.4661646602805707"
$ws.Range("M10").Value = [double]"0.04852733333333333"
$ws.Range("N10").Value = [double]"0.145582"
$ws.Range("O10").Value = [double]"0.0002733940705404138"
$ws.Range("P10").Value = [double]"0.0002733940705404139"
$ws.Range("Q10").Value = [double]"3.019085031940222"
$ws.Range("R10").Value = [double]"27.171765287462"
$ws.Range("S10").Value = [double]"0.0001274466540161944"
$ws.Range("T10").Value = [double]"0.0001274466540161944"

$ws.Range("I11").Value = [double]"0.4661646602805707"
$ws.Range("J11").Value = [double]"0.4661646602805707"
$ws.Range("M11").Value = [double]"73.99751433333334"
$ws.Range("N11").Value = [double]"221.992543"
$ws.Range("O11").Value = [double]"0.4168883856547366"
$ws.Range("P11").Value = [double]"0.4168883856547366"
$ws.Range("Q11").Value = [double]"4603.689767784796"
$ws.Range("R11").Value = [double]"41433.20791006316"
$ws.Range("S11").Value = [double]"0.1943386326736558"
$ws.Range("T11").Value = [double]"0.1943386326736558"

$ws.Range("G12").Value = [double]"7.783044333333334"
$ws.Range("H12").Value = [double]"23.349133"
$ws.Range("I12").Value = [double]"0.05831763893698088"
$ws.Range("J12").Value = [double]"0.05831763893698089"
$ws.Range("M12").Value = [double]"0.1825283333333333"
$ws.Range("N12").Value = [double]"0.547585"
$ws.Range("O12").Value = [double]"0.001028331058213739"
$ws.Range("P12").Value = [double]"0.001028331058213739"
$ws.Range("Q12").Value = [double]"1.420626110422778"
$ws.Range("R12").Value = [double]"12.785634993805"
$ws.Range("S12").Value = [double]"5.996983936059228E-05"
$ws.Range("T12").Value = [double]"5.99698393605923E-05"

$ws.Range("G13").Value = [double]"7.783044333333334"
$ws.Range("H13").Value = [double]"23.349133"
$ws.Range("I13").Value = [double]"0.05831763893698088"
$ws.Range("J13").Value = [double]"0.05831763893698089"
$ws.Range("O13").Value = [double]"0.0001759459539160193"
$ws.Range("P13").Value = [double]"0.0001759459539160193"
$ws.Range("Q13").Value = [double]"0.2430670688781111"
$ws.Range("R13").Value = [double]"2.187603619903"
$ws.Range("S13").Value = [double]"1.026075261289709E-05"
$ws.Range("T13").Value = [double]"1.026075261289709E-05"

$ws.Range("G14").Value = [double]"7.783044333333334"
$ws.Range("H14").Value = [double]"23.349133"
$ws.Range("I14").Value = [double]"0.05831763893698088"
$ws.Range("J14").Value = [double]"0.05831763893698089"
$ws.Range("M14").Value = [double]"103.239782"
$ws.Range("N14").Value = [double]"309.719346"
$ws.Range("O14").Value = [double]"0.5816339432625932"
$ws.Range("P14").Value = [double]"0.5816339432625932"
$ws.Range("Q14").Value = [double]"803.5198002696686"
$ws.Range("R14").Value = [double]"7231.678202427018"
$ws.Range("S14").Value = [double]"0.03391951829668034"
$ws.Range("T14").Value = [double]"0.03391951829668034"

$ws.Range("G15").Value = [double]"7.783044333333334"
$ws.Range("H15").Value = [double]"23.349133"
$ws.Range("I15").Value = [double]"0.05831763893698088"
$ws.Range("J15").Value = [double]"0.05831763893698089"
$ws.Range("M15").Value = [double]"0.04852733333333333"
$ws.Range("N15").Value = [double]"0.145582"
$ws.Range("O15").Value = [double]"0.0002733940705404138"
$ws.Range("P15").Value = [double]"0.0002733940705404139"
$ws.Range("Q15").Value = [double]"0.3776903867117778"
$ws.Range("R15").Value = [double]"3.399213480406"
$ws.Range("S15").Value = [double]"1.594369669328734E-05"
$ws.Range("T15").Value = [double]"1.594369669328734E-05"

$ws.Range("G16").Value = [double]"7.783044333333334"
$ws.Range("H16").Value = [double]"23.349133"
$ws.Range("I16").Value = [double]"0.05831763893698088"
$ws.Range("J16").Value = [double]"0.05831763893698089"
$ws.Range("M16").Value = [double]"73.99751433333334"
$ws.Range("N16").Value = [double]"221.992543"
$ws.Range("O16").Value = [double]"0.4168883856547366"
$ws.Range("P16").Value = [double]"0.4168883856547366"
$ws.Range("Q16").Value = [double]"575.9259346128022"
$ws.Range("R16").Value = [double]"5183.333411515219"
$ws.Range("S16").Value = [double]"0.02431194635163377"
$ws.Range("T16").Value = [double]"0.02431194635163377"

$ws.Range("G17").Value = [double]"30.44016466666666"
$ws.Range("H17").Value = [double]"91.320494"
$ws.Range("I17").Value = [double]"0.2280853681650076"
$ws.Range("J17").Value = [double]"0.2280853681650076"
$ws.Range("M17").Value = [double]"0.1825283333333333"
$ws.Range("N17").Value = [double]"0.547585"
$ws.Range("O17").Value = [double]"0.001028331058213739"
$ws.Range("P17").Value = [double]"0.001028331058213739"
$ws.Range("Q17").Value = [double]"5.556192522998888"
$ws.Range("R17").Value = [double]"50.00573270699"
$ws.Range("S17").Value = [double]"0.0002345472680081925"
$ws.Range("T17").Value = [double]"0.0002345472680081925"

$ws.Range("G18").Value = [double]"30.44016466666666"
$ws.Range("H18").Value = [double]"91.320494"
$ws.Range("I18").Value = [double]"0.2280853681650076"
$ws.Range("J18").Value = [double]"0.2280853681650076"
$ws.Range("O18").Value = [double]"0.0001759459539160193"
$ws.Range("P18").Value = [double]"0.0001759459539160193"
$ws.Range("Q18").Value = [double]"0.9506564892615554"
$ws.Range("R18").Value = [double]"8.555908403354"
$ws.Range("S18").Value = [double]"4.013069767607871E-05"
$ws.Range("T18").Value = [double]"4.013069767607871E-05"

$ws.Range("G19").Value = [double]"30.44016466666666"
$ws.Range("H19").Value = [double]"91.320494"
$ws.Range("I19").Value = [double]"0.2280853681650076"
$ws.Range("J19").Value = [double]"0.2280853681650076"
$ws.Range("M19").Value = [double]"103.239782"
$ws.Range("N19").Value = [double]"309.719346"
$ws.Range("O19").Value = [double]"0.5816339432625932"
$ws.Range("P19").Value = [double]"0.5816339432625932"
$ws.Range("Q19").Value = [double]"3142.635964230769"
$ws.Range("R19").Value = [double]"28283.72367807692"
$ws.Range("S19").Value = [double]"0.1326621920863137"
$ws.Range("T19").Value = [double]"0.1326621920863137"

$ws.Range("G20").Value = [double]"30.44016466666666"
$ws.Range("H20").Value = [double]"91.320494"
$ws.Range("I20").Value = [double]"0.2280853681650076"
$ws.Range("J20").Value = [double]"0.2280853681650076"
$ws.Range("M20").Value = [double]"0.04852733333333333"
$ws.Range("N20").Value = [double]"0.145582"
$ws.Range("O20").Value = [double]"0.0002733940705404138"
$ws.Range("P20").Value = [double]"0.0002733940705404139"
$ws.Range("Q20").Value = [double]"1.477180017500889"
$ws.Range("R20").Value = [double]"13.294620157508"
$ws.Range("S20").Value = [double]"6.235718723334035E-05"
$ws.Range("T20").Value = [double]"6.235718723334036E-05"

$ws.Range("G21").Value = [double]"30.44016466666666"
$ws.Range("H21").Value = [double]"91.320494"
$ws.Range("I21").Value = [double]"0.2280853681650076"
$ws.Range("J21").Value = [double]"0.2280853681650076"
$ws.Range("M21").Value = [double]"73.99751433333334"
$ws.Range("N21").Value = [double]"221.992543"
$ws.Range("O21").Value = [double]"0.4168883856547366"
$ws.Range("P21").Value = [double]"0.4168883856547366"
$ws.Range("Q21").Value = [double]"2252.496521230694"
$ws.Range("R21").Value = [double]"20272.46869107624"
$ws.Range("S21").Value = [double]"0.09508614092577627"
$ws.Range("T21").Value = [double]"0.09508614092577627"

$ws.Range("G22").Value = [double]"10.06688366666667"
$ws.Range("H22").Value = [double]"30.200651"
$ws.Range("I22").Value = [double]"0.07543023805979308"
$ws.Range("J22").Value = [double]"0.07543023805979308"
$ws.Range("M22").Value = [double]"0.1825283333333333"
$ws.Range("N22").Value = [double]"0.547585"
$ws.Range("O22").Value = [double]"0.001028331058213739"
$ws.Range("P22").Value = [double]"0.001028331058213739"
$ws.Range("Q22").Value = [double]"1.837491497537222"
$ws.Range("R22").Value = [double]"16.537423477835"
$ws.Range("S22").Value = [double]"7.756725652534125E-05"
$ws.Range("T22").Value = [double]"7.756725652534126E-05"

$ws.Range("G23").Value = [double]"10.06688366666667"
$ws.Range("H23").Value = [double]"30.200651"
$ws.Range("I23").Value = [double]"0.07543023805979308"
$ws.Range("J23").Value = [double]"0.07543023805979308"
$ws.Range("O23").Value = [double]"0.0001759459539160193"
$ws.Range("P23").Value = [double]"0.0001759459539160193"
$ws.Range("Q23").Value = [double]"0.3143921325378889"
$ws.Range("R23").Value = [double]"2.829529192841"
$ws.Range("S23").Value = [double]"1.327164518954272E-05"
$ws.Range("T23").Value = [double]"1.327164518954272E-05"

$ws.Range("G24").Value = [double]"10.06688366666667"
$ws.Range("H24").Value = [double]"30.200651"
$ws.Range("I24").Value = [double]"0.07543023805979308"
$ws.Range("J24").Value = [double]"0.07543023805979308"
$ws.Range("M24").Value = [double]"103.239782"
$ws.Range("N24").Value = [double]"309.719346"
$ws.Range("O24").Value = [double]"0.5816339432625932"
$ws.Range("P24").Value = [double]"0.5816339432625932"
$ws.Range("Q24").Value = [double]"1039.302875166027"
$ws.Range("R24").Value = [double]"9353.725876494245"
$ws.Range("S24").Value = [double]"0.04387278680395359"
$ws.Range("T24").Value = [double]"0.04387278680395359"

$ws.Range("G25").Value = [double]"10.06688366666667"
$ws.Range("H25").Value = [double]"30.200651"
$ws.Range("I25").Value = [double]"0.07543023805979308"
$ws.Range("J25").Value = [double]"0.07543023805979308"
$ws.Range("M25").Value = [double]"0.04852733333333333"
$ws.Range("N25").Value = [double]"0.145582"
$ws.Range("O25").Value = [double]"0.0002733940705404138"
$ws.Range("P25").Value = [double]"0.0002733940705404139"
$ws.Range("Q25").Value = [double]"0.4885190193202222"
$ws.Range("R25").Value = [double]"4.396671173882"
$ws.Range("S25").Value = [double]"2.062217982499928E-05"
$ws.Range("T25").Value = [double]"2.062217982499928E-05"

$ws.Range("G26").Value = [double]"10.06688366666667"
$ws.Range("H26").Value = [double]"30.200651"
$ws.Range("I26").Value = [double]"0.07543023805979308"
$ws.Range("J26").Value = [double]"0.07543023805979308"
$ws.Range("M26").Value = [double]"73.99751433333334"
$ws.Range("N26").Value = [double]"221.992543"
$ws.Range("O26").Value = [double]"0.4168883856547366"
$ws.Range("P26").Value = [double]"0.4168883856547366"
$ws.Range("Q26").Value = [double]"575.9259346128022"
$ws.Range("R26").Value = [double]"5183.333411515219"
$ws.Range("S26").Value = [double]"0.02431194635163377"
$ws.Range("T26").Value = [double]"0.02431194635163377"

Write-Host "applied updates"